# Updates a set of numeric cells (columns H-N) across the eight craft
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed
# market-board pricing data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2806259.2
$ws.Range("J17").Value = 2871462.8
$ws.Range("L17").Value = 8614388.399999999
$ws.Range("N17").Value = -8614724.399999999
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H96").Value = 581.5417
$ws.Range("I96").Value = 528.625
$ws.Range("K96").Value = 1585.875
$ws.Range("M96").Value = -212.875
$ws.Range("H116").Value = 6965.5
$ws.Range("J116").Value = 6973.75
$ws.Range("L116").Value = 6973.75
$ws.Range("N116").Value = -13857.75
$ws.Range("H132").Value = 4277.278
$ws.Range("I132").Value = 2310.9167
$ws.Range("K132").Value = 6932.750100000001
$ws.Range("M132").Value = -4402.750100000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1235.6666
$ws.Range("I45").Value = 1266
$ws.Range("J45").Value = 1175
$ws.Range("K45").Value = 1266
$ws.Range("L45").Value = 1175
$ws.Range("M45").Value = -889
$ws.Range("N45").Value = -1929
$ws.Range("H61").Value = 6454.9165
$ws.Range("I61").Value = 2522.7144
$ws.Range("J61").Value = 11960
$ws.Range("K61").Value = 2522.7144
$ws.Range("L61").Value = 11960
$ws.Range("M61").Value = -2310.7144
$ws.Range("N61").Value = -12384
$ws.Range("H98").Value = 34701
$ws.Range("J98").Value = 34701
$ws.Range("L98").Value = 34701
$ws.Range("N98").Value = -40691
$ws.Range("H122").Value = 2308.5
$ws.Range("I122").Value = 1817
$ws.Range("J122").Value = 3783
$ws.Range("K122").Value = 5451
$ws.Range("L122").Value = 11349
$ws.Range("M122").Value = -3001
$ws.Range("N122").Value = -16249
$ws.Range("H132").Value = 2478.5
$ws.Range("I132").Value = 2419.9
$ws.Range("J132").Value = 2625
$ws.Range("K132").Value = 7259.700000000001
$ws.Range("L132").Value = 7875
$ws.Range("M132").Value = -4729.700000000001
$ws.Range("N132").Value = -12935
$ws.Range("H136").Value = 6454.9165
$ws.Range("I136").Value = 2522.7144
$ws.Range("J136").Value = 11960
$ws.Range("K136").Value = 7568.1432
$ws.Range("L136").Value = 35880
$ws.Range("M136").Value = -5018.1432
$ws.Range("N136").Value = -40980

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1369.862
$ws.Range("I20").Value = 1223.579
$ws.Range("J20").Value = 1647.8
$ws.Range("K20").Value = 1223.579
$ws.Range("L20").Value = 1647.8
$ws.Range("M20").Value = -976.579
$ws.Range("N20").Value = -2141.8
$ws.Range("H99").Value = 4115.64
$ws.Range("I99").Value = 4115.727
$ws.Range("K99").Value = 4115.727
$ws.Range("M99").Value = -2617.727
$ws.Range("H134").Value = 3361.4
$ws.Range("I134").Value = 2666.8
$ws.Range("J134").Value = 6139.8
$ws.Range("K134").Value = 8000.400000000001
$ws.Range("L134").Value = 18419.4
$ws.Range("M134").Value = -5465.400000000001
$ws.Range("N134").Value = -23489.4

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 437143.12
$ws.Range("I31").Value = 910245.9399999999
$ws.Range("J31").Value = 3465.5833
$ws.Range("K31").Value = 910245.9399999999
$ws.Range("L31").Value = 3465.5833
$ws.Range("M31").Value = -909950.9399999999
$ws.Range("N31").Value = -4055.5833
$ws.Range("H34").Value = 437143.12
$ws.Range("I34").Value = 910245.9399999999
$ws.Range("J34").Value = 3465.5833
$ws.Range("K34").Value = 910245.9399999999
$ws.Range("L34").Value = 3465.5833
$ws.Range("M34").Value = -910043.9399999999
$ws.Range("N34").Value = -3869.5833
$ws.Range("H58").Value = 2601.4285
$ws.Range("I58").Value = 2242.3635
$ws.Range("K58").Value = 2242.3635
$ws.Range("M58").Value = -2039.3635
$ws.Range("H86").Value = 1008407.5
$ws.Range("I86").Value = 1436000.4
$ws.Range("J86").Value = 10690.667
$ws.Range("K86").Value = 1436000.4
$ws.Range("L86").Value = 10690.667
$ws.Range("M86").Value = -1434877.4
$ws.Range("N86").Value = -12936.667
$ws.Range("H89").Value = 1008407.5
$ws.Range("I89").Value = 1436000.4
$ws.Range("J89").Value = 10690.667
$ws.Range("K89").Value = 7180002
$ws.Range("L89").Value = 53453.335
$ws.Range("M89").Value = -7174386
$ws.Range("N89").Value = -64685.335
$ws.Range("H94").Value = 1327.4783
$ws.Range("I94").Value = 810.625
$ws.Range("J94").Value = 1603.1333
$ws.Range("K94").Value = 810.625
$ws.Range("L94").Value = 1603.1333
$ws.Range("M94").Value = -359.625
$ws.Range("N94").Value = -2505.1333
$ws.Range("H122").Value = 3398.111
$ws.Range("I122").Value = 2928.3333
$ws.Range("J122").Value = 4337.6665
$ws.Range("K122").Value = 8784.999899999999
$ws.Range("L122").Value = 13012.9995
$ws.Range("M122").Value = -6334.999899999999
$ws.Range("N122").Value = -17912.9995
$ws.Range("H132").Value = 3523.524
$ws.Range("I132").Value = 3222.5557
$ws.Range("J132").Value = 5329.3335
$ws.Range("K132").Value = 9667.667099999999
$ws.Range("L132").Value = 15988.0005
$ws.Range("M132").Value = -7137.667099999999
$ws.Range("N132").Value = -21048.0005
$ws.Range("H134").Value = 6699.3
$ws.Range("I134").Value = 7219.36
$ws.Range("J134").Value = 4099
$ws.Range("K134").Value = 21658.08
$ws.Range("L134").Value = 12297
$ws.Range("M134").Value = -19123.08
$ws.Range("N134").Value = -17367
$ws.Range("H136").Value = 2601.4285
$ws.Range("I136").Value = 2242.3635
$ws.Range("K136").Value = 6727.0905
$ws.Range("M136").Value = -4177.0905

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 890.25
$ws.Range("I46").Value = 303.14285
$ws.Range("K46").Value = 909.4285500000001
$ws.Range("M46").Value = -818.4285500000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8448.579
$ws.Range("I70").Value = 8382.666999999999
$ws.Range("J70").Value = 8479
$ws.Range("K70").Value = 8382.666999999999
$ws.Range("L70").Value = 8479
$ws.Range("M70").Value = -8112.666999999999
$ws.Range("N70").Value = -9019
$ws.Range("H73").Value = 8448.579
$ws.Range("I73").Value = 8382.666999999999
$ws.Range("J73").Value = 8479
$ws.Range("K73").Value = 8382.666999999999
$ws.Range("L73").Value = 8479
$ws.Range("M73").Value = -7446.666999999999
$ws.Range("N73").Value = -10351
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("N93").Value = 0
$ws.Range("H102").Value = 3575
$ws.Range("I102").Value = 3695.4666
$ws.Range("J102").Value = 3123.25
$ws.Range("K102").Value = 3695.4666
$ws.Range("L102").Value = 3123.25
$ws.Range("M102").Value = -2073.4666
$ws.Range("N102").Value = -6367.25
$ws.Range("H126").Value = 9958.200000000001
$ws.Range("I126").Value = 15704.6
$ws.Range("K126").Value = 47113.8
$ws.Range("M126").Value = -44643.8
$ws.Range("H132").Value = 37667.902
$ws.Range("I132").Value = 45737.707
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 137213.121
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -134683.121
$ws.Range("N132").Value = -35060

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3876.1904
$ws.Range("I68").Value = 4466.6665
$ws.Range("K68").Value = 4466.6665
$ws.Range("M68").Value = -3717.6665
$ws.Range("H71").Value = 3876.1904
$ws.Range("I71").Value = 4466.6665
$ws.Range("K71").Value = 22333.3325
$ws.Range("M71").Value = -18589.3325
$ws.Range("H93").Value = 100002830
$ws.Range("I93").Value = 2793.25
$ws.Range("K93").Value = 2793.25
$ws.Range("M93").Value = -1545.25
$ws.Range("H101").Value = 20748.75
$ws.Range("J101").Value = 20748.75
$ws.Range("L101").Value = 20748.75
$ws.Range("N101").Value = -27238.75
$ws.Range("H122").Value = 4654.778
$ws.Range("I122").Value = 4833
$ws.Range("K122").Value = 14499
$ws.Range("M122").Value = -12049

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 18000
$ws.Range("I5").Value = 50000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 50000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -49888
$ws.Range("N5").Value = -10224
$ws.Range("H31").Value = 14759.5
$ws.Range("J31").Value = 14759.5
$ws.Range("L31").Value = 14759.5
$ws.Range("N31").Value = -15455.5
$ws.Range("H51").Value = 20599.8
$ws.Range("J51").Value = 22999
$ws.Range("L51").Value = 22999
$ws.Range("N51").Value = -24019
$ws.Range("H52").Value = 15599.4
$ws.Range("J52").Value = 22999
$ws.Range("L52").Value = 22999
$ws.Range("N52").Value = -23451
$ws.Range("H61").Value = 37793.145
$ws.Range("J61").Value = 33519
$ws.Range("L61").Value = 33519
$ws.Range("N61").Value = -34103
$ws.Range("H109").Value = 24999.5
$ws.Range("J109").Value = 24999.5
$ws.Range("L109").Value = 24999.5
$ws.Range("N109").Value = -27773.5
